$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 7) to the student diagnostics sheet
$rowValues = @(
    6,
    "Masculino",
    20,
    "Costa",
    "Urbana",
    "VI",
    "B",
    "tarde",
    37,
    "alto",
    "A menudo",
    "De vez en cuando",
    "A menudo",
    "Casi nunca",
    "De vez en cuando",
    "Casi nunca",
    "De vez en cuando",
    "De vez en cuando",
    "De vez en cuando",
    "De vez en cuando",
    "A menudo",
    "Casi nunca Casi nunca",
    "De vez en cuando",
    "De vez en cuando"
)

$targetRow = 7
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item($targetRow, $i + 1).Value = $rowValues[$i]
}

Write-Host ("New UsedRange: " + $ws.UsedRange.Address())
